$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 6).Value2 = 2.04
$ws.Cells.Item(2, 7).Value2 = 2.22
$ws.Cells.Item(2, 8).Value2 = 3.15
$ws.Cells.Item(2, 9).Value2 = 3.65
$ws.Cells.Item(2, 10).Value2 = 3.95
$ws.Cells.Item(2, 11).Value2 = 4.6
$ws.Cells.Item(2, 14).Value2 = 5.8
$ws.Cells.Item(2, 15).Value2 = 1.16
$ws.Cells.Item(2, 16).Value2 = 2.68
$ws.Cells.Item(2, 17).Value2 = 1.48
$ws.Cells.Item(2, 18).Value2 = 1.7
$ws.Cells.Item(2, 19).Value2 = 2.22
$ws.Cells.Item(2, 20).Value2 = 1.49
$ws.Cells.Item(2, 21).Value2 = 2.6
$ws.Cells.Item(2, 22).Value2 = 1.38
$ws.Cells.Item(2, 23).Value2 = 1.83
$ws.Cells.Item(2, 25).Value2 = 22
$ws.Cells.Item(2, 26).Value2 = 32
$ws.Cells.Item(2, 28).Value2 = 16.5
$ws.Cells.Item(2, 29).Value2 = 11
$ws.Cells.Item(2, 30).Value2 = 16
$ws.Cells.Item(2, 31).Value2 = 85
$ws.Cells.Item(2, 32).Value2 = 18.5
$ws.Cells.Item(2, 33).Value2 = 12.5
$ws.Cells.Item(2, 34).Value2 = 15.5
$ws.Cells.Item(2, 35).Value2 = 85
$ws.Cells.Item(2, 36).Value2 = 36
$ws.Cells.Item(2, 37).Value2 = 20
$ws.Cells.Item(2, 38).Value2 = 36
$ws.Cells.Item(2, 39).Value2 = 580
$ws.Cells.Item(2, 40).Value2 = 9.6
$ws.Cells.Item(2, 41).Value2 = 21

# Row 3
$ws.Cells.Item(3, 14).Value2 = 1.11

# Row 4
$ws.Cells.Item(4, 6).Value2 = 2.6
$ws.Cells.Item(4, 7).Value2 = 2.72
$ws.Cells.Item(4, 8).Value2 = 3
$ws.Cells.Item(4, 9).Value2 = 3.15
$ws.Cells.Item(4, 10).Value2 = 3.25
$ws.Cells.Item(4, 12).Value2 = 1.51
$ws.Cells.Item(4, 13).Value2 = 1.09
$ws.Cells.Item(4, 14).Value2 = 3.1
$ws.Cells.Item(4, 15).Value2 = 1.43
$ws.Cells.Item(4, 16).Value2 = 1.69
$ws.Cells.Item(4, 19).Value2 = 4.5
$ws.Cells.Item(4, 20).Value2 = 1.88
$ws.Cells.Item(4, 21).Value2 = 1.94
$ws.Cells.Item(4, 22).Value2 = 1.46
$ws.Cells.Item(4, 23).Value2 = 1.58
$ws.Cells.Item(4, 27).Value2 = 340
$ws.Cells.Item(4, 32).Value2 = 17.5
$ws.Cells.Item(4, 35).Value2 = 200
$ws.Cells.Item(4, 40).Value2 = 42

# Row 5
$ws.Cells.Item(5, 7).Value2 = 1.95
$ws.Cells.Item(5, 8).Value2 = 5.3
$ws.Cells.Item(5, 10).Value2 = 3.2
$ws.Cells.Item(5, 11).Value2 = 3.5
$ws.Cells.Item(5, 12).Value2 = 1.59
$ws.Cells.Item(5, 14).Value2 = 2.66
$ws.Cells.Item(5, 16).Value2 = 1.54
$ws.Cells.Item(5, 17).Value2 = 2.68
$ws.Cells.Item(5, 18).Value2 = 1.19
$ws.Cells.Item(5, 19).Value2 = 5.5
$ws.Cells.Item(5, 21).Value2 = 1.66
$ws.Cells.Item(5, 24).Value2 = 9.4
$ws.Cells.Item(5, 25).Value2 = 26
$ws.Cells.Item(5, 27).Value2 = 220
$ws.Cells.Item(5, 29).Value2 = 8.199999999999999
$ws.Cells.Item(5, 31).Value2 = 140
$ws.Cells.Item(5, 32).Value2 = 10.5
$ws.Cells.Item(5, 33).Value2 = 11.5
$ws.Cells.Item(5, 36).Value2 = 130
$ws.Cells.Item(5, 37).Value2 = 980
$ws.Cells.Item(5, 39).Value2 = 980

# Row 6
$ws.Cells.Item(6, 6).Value2 = 1.88
$ws.Cells.Item(6, 8).Value2 = 4.9
$ws.Cells.Item(6, 10).Value2 = 3.45
$ws.Cells.Item(6, 11).Value2 = 3.7
$ws.Cells.Item(6, 12).Value2 = 1.47
$ws.Cells.Item(6, 14).Value2 = 3.25
$ws.Cells.Item(6, 16).Value2 = 1.75
$ws.Cells.Item(6, 18).Value2 = 1.28
$ws.Cells.Item(6, 19).Value2 = 4.2
$ws.Cells.Item(6, 20).Value2 = 2.04
$ws.Cells.Item(6, 21).Value2 = 1.87
$ws.Cells.Item(6, 24).Value2 = 13.5
$ws.Cells.Item(6, 25).Value2 = 17
$ws.Cells.Item(6, 26).Value2 = 55
$ws.Cells.Item(6, 27).Value2 = 170
$ws.Cells.Item(6, 28).Value2 = 8.4
$ws.Cells.Item(6, 29).Value2 = 8.6
$ws.Cells.Item(6, 30).Value2 = 24
$ws.Cells.Item(6, 31).Value2 = 1000
$ws.Cells.Item(6, 32).Value2 = 11.5
$ws.Cells.Item(6, 33).Value2 = 10.5
$ws.Cells.Item(6, 36).Value2 = 22
$ws.Cells.Item(6, 37).Value2 = 26
$ws.Cells.Item(6, 40).Value2 = 16

# Row 7
$ws.Cells.Item(7, 6).Value2 = 2.98
$ws.Cells.Item(7, 7).Value2 = 3.4
$ws.Cells.Item(7, 9).Value2 = 2.62
$ws.Cells.Item(7, 10).Value2 = 3.35
$ws.Cells.Item(7, 12).Value2 = 1.44
$ws.Cells.Item(7, 13).Value2 = 1.08
$ws.Cells.Item(7, 14).Value2 = 3.65
$ws.Cells.Item(7, 15).Value2 = 1.34
$ws.Cells.Item(7, 16).Value2 = 1.91
$ws.Cells.Item(7, 17).Value2 = 2.02
$ws.Cells.Item(7, 19).Value2 = 3.65
$ws.Cells.Item(7, 20).Value2 = 1.81
$ws.Cells.Item(7, 22).Value2 = 1.63
$ws.Cells.Item(7, 23).Value2 = 1.42
$ws.Cells.Item(7, 24).Value2 = 15
$ws.Cells.Item(7, 27).Value2 = 36
$ws.Cells.Item(7, 28).Value2 = 12.5
$ws.Cells.Item(7, 29).Value2 = 8.4
$ws.Cells.Item(7, 30).Value2 = 13.5
$ws.Cells.Item(7, 31).Value2 = 980
$ws.Cells.Item(7, 33).Value2 = 14.5
$ws.Cells.Item(7, 36).Value2 = 60
$ws.Cells.Item(7, 37).Value2 = 42
$ws.Cells.Item(7, 39).Value2 = 140
$ws.Cells.Item(7, 40).Value2 = 48
$ws.Cells.Item(7, 41).Value2 = 24

# Row 8
$ws.Cells.Item(8, 7).Value2 = 4
$ws.Cells.Item(8, 9).Value2 = 2.38
$ws.Cells.Item(8, 14).Value2 = 2.84
$ws.Cells.Item(8, 16).Value2 = 1.62
$ws.Cells.Item(8, 17).Value2 = 2.4
$ws.Cells.Item(8, 21).Value2 = 1.92
$ws.Cells.Item(8, 22).Value2 = 1.72
$ws.Cells.Item(8, 23).Value2 = 1.33
$ws.Cells.Item(8, 24).Value2 = 10
$ws.Cells.Item(8, 25).Value2 = 8.199999999999999
$ws.Cells.Item(8, 31).Value2 = 1000
$ws.Cells.Item(8, 33).Value2 = 19.5
$ws.Cells.Item(8, 35).Value2 = 55
$ws.Cells.Item(8, 39).Value2 = 200

# Row 9
$ws.Cells.Item(9, 6).Value2 = 3.2
$ws.Cells.Item(9, 7).Value2 = 4.2
$ws.Cells.Item(9, 8).Value2 = 2.36
$ws.Cells.Item(9, 9).Value2 = 2.48
$ws.Cells.Item(9, 10).Value2 = 2.88
$ws.Cells.Item(9, 11).Value2 = 3.5
$ws.Cells.Item(9, 12).Value2 = 1.46
$ws.Cells.Item(9, 13).Value2 = 1.09
$ws.Cells.Item(9, 14).Value2 = 2.74
$ws.Cells.Item(9, 15).Value2 = 1.39
$ws.Cells.Item(9, 16).Value2 = 1.77
$ws.Cells.Item(9, 17).Value2 = 2.08
$ws.Cells.Item(9, 18).Value2 = 1.24
$ws.Cells.Item(9, 19).Value2 = 3.45
$ws.Cells.Item(9, 20).Value2 = 1.89
$ws.Cells.Item(9, 21).Value2 = 1.9
$ws.Cells.Item(9, 22).Value2 = 1.67
$ws.Cells.Item(9, 23).Value2 = 1.32
